# Apply the "feat: add 2022-Q4 data" edit:
#  1. Insert a brand-new worksheet "2022-Q4" right after "总计", containing the
#     16-fund holdings detail for that quarter (same layout/style as the other
#     quarterly sheets, e.g. "2022-Q3").
#  2. Update the "总计" (summary) sheet so it gets a new first data row for
#     2022-Q4 and all the existing quarters shift down one row (their index
#     in column A, and the physical row number, both increment by one), with
#     "2021-Q1" reappearing (it had fallen out of the old summary) right
#     before "2020-Q4".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "2022-Q4" worksheet, placed immediately after "总计" (first tab).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $totalSheet)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q4.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# Bold + thin box border + centered/top aligned header row (matches the
# other quarter sheets' header formatting).
$headerRange = $q4.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Fund rows for 2022-Q4 (code, name, scale, stock position, weight%, held
# value (CNY 100M), position rank). Columns C-G are kept as text to match
# the other quarterly sheets' storage (inline/shared strings there too).
$rows = @(
    @("213001", "宝盈鸿利收益灵活配置混合A",          "14.47", "91.18", "4.34", "0.6280", 7),
    @("000986", "太平灵活配置混合",                    "11.54", "81.01", "4.61", "0.5320", 8),
    @("010328", "博时荣华灵活配置混合A",                "4.01",  "72.54", "2.36", "0.0946", 9),
    @("005933", "新疆前海联合先进制造灵活配置混合A",      "0.90",  "90.73", "4.16", "0.0374", 10),
    @("000066", "诺安鸿鑫混合A",                        "0.70",  "90.20", "5.01", "0.0351", 4),
    @("006429", "诺安恒鑫混合",                         "0.92",  "66.36", "3.60", "0.0331", 8),
    @("009537", "太平行业优选股票A",                     "0.55",  "93.57", "5.31", "0.0292", 8),
    @("007581", "宝盈鸿利收益灵活配置混合C",             "0.62",  "91.18", "4.34", "0.0269", 7),
    @("009538", "太平行业优选股票C",                     "0.18",  "93.57", "5.31", "0.0096", 8),
    @("008629", "大成景瑞稳健配置混合A",                 "0.66",  "29.66", "1.38", "0.0091", 10),
    @("014608", "中欧周期景气混合A",                     "0.20",  "93.00", "3.91", "0.0078", 7),
    @("005934", "新疆前海联合先进制造灵活配置混合C",       "0.09",  "90.73", "4.16", "0.0037", 10),
    @("008630", "大成景瑞稳健配置混合C",                 "0.25",  "29.66", "1.38", "0.0034", 10),
    @("010329", "博时荣华灵活配置混合C",                 "0.14",  "72.54", "2.36", "0.0033", 9),
    @("014609", "中欧周期景气混合C",                     "0.05",  "93.00", "3.91", "0.0020", 7),
    @("014498", "诺安鸿鑫混合C",                        "0.01",  "90.20", "5.01", "0.0005", 4)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    # Column A: 0-based running index, bold+bordered like the other sheets.
    $q4.Cells.Item($r, 1).Value = $i
    $aCell = $q4.Range($q4.Cells.Item($r, 1), $q4.Cells.Item($r, 1))
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    $q4.Cells.Item($r, 2).Value = $data[0]

    $q4.Cells.Item($r, 3).Value = $data[1]

    $q4.Cells.Item($r, 4).NumberFormat = "@"
    $q4.Cells.Item($r, 4).Value = $data[2]

    $q4.Cells.Item($r, 5).NumberFormat = "@"
    $q4.Cells.Item($r, 5).Value = $data[3]

    $q4.Cells.Item($r, 6).NumberFormat = "@"
    $q4.Cells.Item($r, 6).Value = $data[4]

    $q4.Cells.Item($r, 7).NumberFormat = "@"
    $q4.Cells.Item($r, 7).Value = $data[5]

    $q4.Cells.Item($r, 8).Value = $data[6]
}

# ---------------------------------------------------------------------------
# 2. Rewrite "总计" summary rows 2-10: 2022-Q4 is the new top row, every
#    other quarter shifts down by one, and 2021-Q1 (previously dropped)
#    reappears just above 2020-Q4.
# ---------------------------------------------------------------------------
$summary = $totalSheet

$summaryRows = @(
    @("2022-Q4", 16, 1.46),
    @("2022-Q3", 34, 4.85),
    @("2022-Q2", 19, 3.33),
    @("2022-Q1", 6,  1.38),
    @("2021-Q4", 17, 10.69),
    @("2021-Q3", 64, 34.13),
    @("2021-Q2", 35, 6.62),
    @("2021-Q1", 15, 3.94),
    @("2020-Q4", 1,  0.01)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $data = $summaryRows[$i]

    $summary.Cells.Item($r, 1).Value = $i
    $aCell = $summary.Range($summary.Cells.Item($r, 1), $summary.Cells.Item($r, 1))
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    $summary.Cells.Item($r, 2).Value = $data[0]
    $summary.Cells.Item($r, 3).Value = $data[1]
    $summary.Cells.Item($r, 4).Value = $data[2]
}
